# Updates cryptos list figures (price + 1h volume change) to the latest
# scrape, and swaps the FraxShare/MXToken ranking rows, matching the
# "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.560.65"
$ws.Range("E2").Value = "  +0.80%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.019.11"
$ws.Range("E3").Value = "  +0.85%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.12%  "

# --- Row 5: BNB ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "262.94"
$ws.Range("E5").Value = "  +6.39%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  -1.98%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  -0.05%  "

# --- Row 8: Solana ---
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.11"
$ws.Range("E8").Value = "  -6.71%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +0.65%  "

# --- Row 10: Dogecoin ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -3.55%  "

# --- Row 12: Chainlink ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.39"
$ws.Range("E12").Value = "  -3.88%  "

# --- Row 13: WrappedliquidstakedEther2.0 ---
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.313.76"
$ws.Range("E13").Value = "  +0.65%  "

# --- Row 14: Polygon ---
$ws.Range("E14").Value = "  -4.75%  "

# --- Row 15: Avalanche ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.85"
$ws.Range("E15").Value = "  -8.17%  "

# --- Row 16: Polkadot ---
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  -4.01%  "

# --- Row 17: WrappedEther ---
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.012.09"
$ws.Range("E17").Value = "  +0.05%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.471.93"
$ws.Range("E18").Value = "  +0.87%  "

# --- Row 19: Litecoin ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.67"
$ws.Range("E19").Value = "  -1.03%  "

# --- Row 20: ShibaInu ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0842"
$ws.Range("E20").Value = "  -2.61%  "

# --- Row 21: Uniswap ---
$ws.Range("E21").Value = "  -0.42%  "

# --- Row 22: BitcoinCash ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.78"
$ws.Range("E22").Value = "  -0.84%  "

# --- Row 23: PancakeSwap ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.69"
$ws.Range("E23").Value = "  +8.22%  "

# --- Row 24: Dai ---
$ws.Range("E24").Value = "  -0.03%  "

# --- Row 25: Toncoin ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  -1.48%  "

# --- Row 26: Monero ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.50"
$ws.Range("E26").Value = "  +0.60%  "

# --- Row 27: Cosmos ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.94"
$ws.Range("E27").Value = "  -5.00%  "

# --- Row 28: EthereumClassic ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.71"
$ws.Range("E28").Value = "  +0.14%  "

# --- Row 29: Kaspa ---
$ws.Range("E29").Value = "  -10.96%  "

# --- Row 30: ImmutableX ---
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.32"
$ws.Range("E30").Value = "  -1.43%  "

# --- Row 31: Stellar ---
$ws.Range("E31").Value = "  -1.16%  "

# --- Row 32: Hedera ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0652"
$ws.Range("E32").Value = "  -1.29%  "

# --- Row 33: Filecoin ---
$ws.Range("E33").Value = "  -3.66%  "

# --- Row 34: InternetComputer(DFINITY) ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").Value = "  +0.36%  "

# --- Row 35: LidoDAOToken ---
$ws.Range("E35").Value = "  +1.67%  "

# --- Row 36: WEMIXToken ---
$ws.Range("E36").Value = "  +1.17%  "

# --- Row 37: BinanceUSD ---
$ws.Range("E37").Value = "  -0.12%  "

# --- Row 38: RenderToken ---
$ws.Range("E38").Value = "  +1.95%  "

# --- Row 39: THORChain ---
$ws.Range("E39").Value = "  -5.25%  "

# --- Row 40: HuobiToken ---
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.06"
$ws.Range("E40").Value = "  +4.51%  "

# --- Row 41: TrustWalletToken ---
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.21"
$ws.Range("E41").Value = "  +2.68%  "

# --- Row 42: Cronos ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0936"
$ws.Range("E42").Value = "  -4.46%  "

# --- Row 43: VeChain ---
$ws.Range("E43").Value = "  -1.03%  "

# --- Row 44: Maker ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.392.46"
$ws.Range("E44").Value = "  +1.08%  "

# --- Row 45: Aave ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.08"
$ws.Range("E45").Value = "  -1.19%  "

# --- Row 46: InjectiveProtocol ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.71"
$ws.Range("E46").Value = "  -6.40%  "

# --- Row 47: ARBITRUM ---
$ws.Range("E47").Value = "  -1.84%  "

# --- Row 48/49: FraxShare <-> MXToken swap ranking, with refreshed figures ---
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.93"
$ws.Range("E48").Value = "  +2.70%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.05"
$ws.Range("E49").Value = "  -3.07%  "

# --- Row 50: RocketPoolETH ---
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.205.70"
$ws.Range("E50").Value = "  +0.62%  "

# --- Row 51: NEARProtocol ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.96"
$ws.Range("E51").Value = "  -3.65%  "
